$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3587
$ws1.Range("F5").Value = 2215
$ws1.Range("F9").Value = 78
$ws1.Range("F10").Value = 67
$ws1.Range("F13").Value = 1897

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3587
$ws4.Range("F5").Value = 2215
$ws4.Range("F10").Value = 78
$ws4.Range("F11").Value = 67
$ws4.Range("F16").Value = 1897
